{"js": "// Apply the four wording/phrasing corrections to the \"career resume\"\n// document. Each target paragraph is held entirely within a single run,\n// so we locate each old sentence with a body search and swap it for the\n// corrected sentence in place (preserving the run's formatting).\n\nconst replacements = [\n  {\n    oldText:\n      \"The working utility invention; Natural Human WaveLength & Impedance Meter/Visual Recognition Medical Instrument I formulated, designed, engineered and manufacture, submitted a white paper to U.S Army Research Laboratory(ARL).\",\n    newText:\n      \"The utility invention; Natural Human WaveLength & Impedance Meter/Visual Recognition Medical Instrument I formulated, designed, engineered and industriallyl manufactured, submitted a white paper to U.S Army Research Laboratory(ARL).\",\n  },\n  {\n    oldText:\n      \"First and only to achieve Nuclear Fusion by implicitly formulating, designing, engineering and manufacture the working utility invention; Nuclear Fusion Reactor and Fusion Reactionary Engine, communicated to U.S. National Aeronautical & Space Agency(NASA) for geo-satellite concatenation of scientific proofs. \",\n    newText:\n      \"First and only to achieve Nuclear Fusion by implicitly formulating, designing, engineering and industrially manufacturing the utility invention the Nuclear Fusion Reactor and the Fusion Reactionary Engine; communicated to U.S. National Aeronautical & Space Agency(NASA) for geo-satellite concatenation of scientific proofs. \",\n  },\n  {\n    oldText:\n      \"Formulated, designing,  engineered and partially manufacture the working utility invention; Full-Spatial Median-Free Liquid and Photonic Bit Transfer Module System 276,480-bit Computational Processor/Computer and submitted working draft on its architecture to U.S. National Security Agency(NSA).\",\n    newText:\n      \"Formulated, designing,  engineered and industrially manufacturing ready the utility invention; Full-Spatial Median-Free Liquid and Photonic Bit Transfer Module System 276,480-bit Computational Processor/Computer and submitted working draft on its architecture to U.S. National Security Agency(NSA).\",\n  },\n  {\n    oldText:\n      \"Technical report to U.S. Naval Research Laboratory(NRL) on the F-23 jet fighter aircrafts LIDAR sub-systems engineering flaw.\",\n    newText:\n      \"Technical report to U.S. Naval Research Laboratory(NRL) on diagnosed F-23n jet fighter aircrafts; LIDAR sub-system engineering flaw.\",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { oldText, newText } of replacements) {\n  const searchResults = body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  if (searchResults.items.length === 0) {\n    throw new Error(\"Could not find expected text: \" + oldText);\n  }\n\n  for (const range of searchResults.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the four wording/phrasing corrections to the \"career resume\"\n# document. Each target paragraph is held entirely within a single run,\n# so Find/Replace against the full old sentence and swap it for the\n# corrected sentence (formatting of the run is left untouched).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{\n        Old = \"The working utility invention; Natural Human WaveLength & Impedance Meter/Visual Recognition Medical Instrument I formulated, designed, engineered and manufacture, submitted a white paper to U.S Army Research Laboratory(ARL).\"\n        New = \"The utility invention; Natural Human WaveLength & Impedance Meter/Visual Recognition Medical Instrument I formulated, designed, engineered and industriallyl manufactured, submitted a white paper to U.S Army Research Laboratory(ARL).\"\n    },\n    @{\n        Old = \"First and only to achieve Nuclear Fusion by implicitly formulating, designing, engineering and manufacture the working utility invention; Nuclear Fusion Reactor and Fusion Reactionary Engine, communicated to U.S. National Aeronautical & Space Agency(NASA) for geo-satellite concatenation of scientific proofs. \"\n        New = \"First and only to achieve Nuclear Fusion by implicitly formulating, designing, engineering and industrially manufacturing the utility invention the Nuclear Fusion Reactor and the Fusion Reactionary Engine; communicated to U.S. National Aeronautical & Space Agency(NASA) for geo-satellite concatenation of scientific proofs. \"\n    },\n    @{\n        Old = \"Formulated, designing,  engineered and partially manufacture the working utility invention; Full-Spatial Median-Free Liquid and Photonic Bit Transfer Module System 276,480-bit Computational Processor/Computer and submitted working draft on its architecture to U.S. National Security Agency(NSA).\"\n        New = \"Formulated, designing,  engineered and industrially manufacturing ready the utility invention; Full-Spatial Median-Free Liquid and Photonic Bit Transfer Module System 276,480-bit Computational Processor/Computer and submitted working draft on its architecture to U.S. National Security Agency(NSA).\"\n    },\n    @{\n        Old = \"Technical report to U.S. Naval Research Laboratory(NRL) on the F-23 jet fighter aircrafts LIDAR sub-systems engineering flaw.\"\n        New = \"Technical report to U.S. Naval Research Laboratory(NRL) on diagnosed F-23n jet fighter aircrafts; LIDAR sub-system engineering flaw.\"\n    }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
